{"js": "// Reposition the floating picture and insert a new empty paragraph\n// right after it (before the \"\u041e\u043f\u0438\u0441\u0430\u043d\u0438\u0435 \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u0435\u043c\u043e\u0433\u043e \u043e\u0431\u043e\u0440\u0443\u0434\u043e\u0432\u0430\u043d\u0438\u044f:\"\n// paragraph), matching the target OOXML diff:\n//   - <wp:positionH><wp:align>center</wp:align></wp:positionH>\n//       -> <wp:positionH><wp:posOffset>0</wp:posOffset></wp:positionH>\n//   - <wp:positionV><wp:posOffset>635</wp:posOffset></wp:positionV>\n//       -> <wp:positionV><wp:posOffset>-720090</wp:posOffset></wp:positionV>\n//   - a new empty paragraph (Normal style, empty run) inserted after the\n//     paragraph that holds the picture.\n\nconst body = context.document.body;\n\n// --- 1. Reposition the floating shape (picture) ---------------------------\nconst shapes = body.shapes;\nshapes.load(\"items\");\nawait context.sync();\n\nif (shapes.items.length > 0) {\n  const pic = shapes.items[0];\n  // 1 pt = 12700 EMU, so 0 EMU -> 0 pt, -720090 EMU -> -56.7 pt.\n  pic.left = 0;\n  pic.top = -56.7;\n  await context.sync();\n}\n\n// --- 2. Insert a new empty paragraph after the picture's paragraph --------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst picParagraph = paragraphs.items[0];\npicParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n", "ps1": "# Reposition the floating picture and insert a new empty paragraph right\n# after it (before the \"\u041e\u043f\u0438\u0441\u0430\u043d\u0438\u0435 \u0438\u0441\u043f\u043e\u043b\u044c\u0437\u0443\u0435\u043c\u043e\u0433\u043e \u043e\u0431\u043e\u0440\u0443\u0434\u043e\u0432\u0430\u043d\u0438\u044f:\" paragraph),\n# matching the target OOXML diff:\n#   - <wp:positionH><wp:align>center</wp:align></wp:positionH>\n#       -> <wp:positionH><wp:posOffset>0</wp:posOffset></wp:positionH>\n#   - <wp:positionV><wp:posOffset>635</wp:posOffset></wp:positionV>\n#       -> <wp:positionV><wp:posOffset>-720090</wp:posOffset></wp:positionV>\n#   - a new empty paragraph (Normal style, empty run) inserted after the\n#     paragraph that holds the picture.\n\n$d = $word.ActiveDocument\n\n# --- 1. Reposition the floating shape (picture) ----------------------------\n$shp = $d.Shapes(1)\n# 1 pt = 12700 EMU, so 0 EMU -> 0 pt, -720090 EMU -> -56.7 pt.\n$shp.Left = 0\n$shp.Top = -56.7\n\n# --- 2. Insert a new empty paragraph after the picture's paragraph --------\n$picParagraph = $d.Paragraphs(1)\n$picParagraph.Range.InsertParagraphAfter()\n"}
